# Agregué una diapositiva mas
# Otra diapositiva

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the "datetimeFigureOut" date placeholder (8/23/13 -> 10/17/14)
#    on the slide master and on every slide layout.
# ---------------------------------------------------------------------
function Update-DatePlaceholders {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes "10/17/14"

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholders $layout.Shapes "10/17/14"
}

# ---------------------------------------------------------------------
# 2) Slide 1 ("Ingenier" + "ía de Software I" -> single run)
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleRange = $slide1.Shapes.Item(1).TextFrame.TextRange
# Force a real content change first so the run collapses into one run
# instead of being left untouched because the text already matches.
$titleRange.Text = "placeholder"
$titleRange = $slide1.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Ingeniería de Software I`t"

# ---------------------------------------------------------------------
# 3) Slide 2 content placeholder paragraph run merges
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$bodyRange = $slide2.Shapes.Item(2).TextFrame.TextRange

$para2 = $bodyRange.Paragraphs(2, 1)
$para2.Text = "placeholder"
$para2 = $bodyRange.Paragraphs(2, 1)
$para2.Text = "¿Qué es la IS y para que sirve?"

$para3 = $bodyRange.Paragraphs(3, 1)
$para3.Text = "placeholder"
$para3 = $bodyRange.Paragraphs(3, 1)
$para3.Text = "Mencionar las características de un buen software"

$para5 = $bodyRange.Paragraphs(5, 1)
$para5.Text = "placeholder"
$para5 = $bodyRange.Paragraphs(5, 1)
$para5.Text = "Mencione y describa por lo menos 5 mejores prácticas que ayuden a mitigar mitos en la IS"

# ---------------------------------------------------------------------
# 4) Add the new third slide ("Otra diapositiva" / "otra diapositiva")
#    using the same "Título y objetos" layout as slide 2.
# ---------------------------------------------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

$newTitle = $newSlide.Shapes.Item(1).TextFrame.TextRange
$newTitle.Text = "Otra diapositiva"
$newTitle.LanguageID = "es-ES"

$newBody = $newSlide.Shapes.Item(2).TextFrame.TextRange
$newBody.Text = "otra diapositiva"
$newBody.LanguageID = "es-ES"

Write-Output "done"
